# Append the "2021年" row (row 11) to Sheet1, mirroring the existing
# year rows above it (A2:AQ10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: copy the style from the row above (bold / centered / bordered
# "year" label style) and then overwrite the value with the new year.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = "2021年"

# Column E has no data point for this year (matches the existing blank
# cells in E6:E10) — write it as a genuine empty string value rather than
# leaving the cell completely absent.
$ws.Range("E11").Value = "'"
$ws.Range("E11").ClearFormats()

# Remaining columns: plain numeric data for 2021.
$ws.Range("B11").Value = 1842.52
$ws.Range("C11").Value = 283.95
$ws.Range("D11").Value = 178.91
$ws.Range("F11").Value = 543.99
$ws.Range("G11").Value = 1665.02
$ws.Range("H11").Value = 155.09
$ws.Range("I11").Value = 697.03
$ws.Range("J11").Value = 104.05
$ws.Range("K11").Value = 38356.15
$ws.Range("L11").Value = 7.96
$ws.Range("M11").Value = 57.23
$ws.Range("N11").Value = 60.71
$ws.Range("O11").Value = 148.2
$ws.Range("P11").Value = 2789.48
$ws.Range("Q11").Value = 143.93
$ws.Range("R11").Value = 38.9
$ws.Range("S11").Value = 224.1
$ws.Range("T11").Value = 206.81
$ws.Range("U11").Value = 3155.99
$ws.Range("V11").Value = 3990.53
$ws.Range("W11").Value = 667.97
$ws.Range("X11").Value = 179.1
$ws.Range("Y11").Value = 1792.74
$ws.Range("Z11").Value = 1320.46
$ws.Range("AA11").Value = 10.8
$ws.Range("AB11").Value = 2324.14
$ws.Range("AC11").Value = 139.77
$ws.Range("AD11").Value = 149.13
$ws.Range("AE11").Value = 42.22
$ws.Range("AF11").Value = 3328.34
$ws.Range("AG11").Value = 1641.55
$ws.Range("AH11").Value = 154.97
$ws.Range("AI11").Value = 1393.78
$ws.Range("AJ11").Value = 208.33
$ws.Range("AK11").Value = 865.12
$ws.Range("AL11").Value = 3699.63
$ws.Range("AM11").Value = 814.67
$ws.Range("AN11").Value = 114.74
$ws.Range("AO11").Value = 162.69
$ws.Range("AP11").Value = 2965.87
$ws.Range("AQ11").Value = 84.88
